$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 557 (existing data 557:576 shifts
# down to 559:578, matching the original sheet's column formatting/styles).
$ws.Rows("557:558").Insert()

# Row 557: new weekly entry ("$/caja 36 atados" unit line)
$ws.Range("A557").Value = 6
$ws.Range("B557").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C557").Value = "Metropolitana"
$ws.Range("D557").Value = 44509
$ws.Range("E557").Value = 13
$ws.Range("F557").Value = 100112040
$ws.Range("G557").Value = "Cilantro"
$ws.Range("H557").Value = "Sin especificar"
$ws.Range("I557").Value = "Primera"
$ws.Range("J557").Value = 650
$ws.Range("K557").Value = 4500
$ws.Range("L557").Value = 5000
$ws.Range("M557").Value = 4715
$ws.Range("N557").Value = "$/caja 36 atados"
$ws.Range("O557").Value = "Región Metropolitana"
$ws.Range("P557").Value = 131
$ws.Range("Q557").Value = 36
$ws.Range("R557").Value = "Hortaliza"

# Row 558: new weekly entry ("$/docena de atados" unit line)
$ws.Range("A558").Value = 6
$ws.Range("B558").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C558").Value = "Metropolitana"
$ws.Range("D558").Value = 44509
$ws.Range("E558").Value = 13
$ws.Range("F558").Value = 100112040
$ws.Range("G558").Value = "Cilantro"
$ws.Range("H558").Value = "Sin especificar"
$ws.Range("I558").Value = "Primera"
$ws.Range("J558").Value = 340
$ws.Range("K558").Value = 11000
$ws.Range("L558").Value = 12000
$ws.Range("M558").Value = 11441
$ws.Range("N558").Value = "$/docena de atados"
$ws.Range("O558").Value = "Región Metropolitana"
$ws.Range("P558").Value = 3814
$ws.Range("Q558").Value = 3
$ws.Range("R558").Value = "Hortaliza"
